# The test fixture's "Units" column (L) held a leftover
# "1-4 Parlier field response score" value that isn't used anywhere else in
# the workbook; clear it out of L2:L4 (this also drops the now-orphaned
# shared string from the saved file's string table).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("L2:L4").ClearContents()

# Leave the sheet scrolled/selected the way the author had it when saving.
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("L2:L4").Select()
